# Insert a new medicine row ("FUCIDIN 2% CREAM 30 GM") into the stock report.
# The new row goes right after row 7 (DOLCYL M ...) and before the old row 8
# (GLYBOFEN ...), i.e. it becomes the new row 8; everything below shifts down
# by one row, and the totals row is recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before row 8 (old row 8 "GLYBOFEN..." and everything
# after it shifts down by one row).
$ws.Rows.Item(8).Insert()

# The inserted row doesn't automatically pick up the bordered data-row
# format, so copy it explicitly from row 7 (same look as every other data
# row: A=style6, B:G=style7, H:K=style8, L:M=style9, N=style7).
$ws.Range("A7:N7").Copy()
$ws.Range("A8:N8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(8).RowHeight = $ws.Rows.Item(7).RowHeight

# Populate the new row 8 with the new medicine's data.
$ws.Cells.Item(8, 1).Value2 = 5          # "م" sequence number
$ws.Cells.Item(8, 2).Value2 = "FUCIDIN 2% CREAM 30 GM"
$ws.Cells.Item(8, 8).Value2 = "1:0"
$ws.Cells.Item(8, 12).Value2 = -96
$ws.Cells.Item(8, 14).Value2 = "1:0"

# Recreate the merged cell ranges for the new row (mirrors the other rows).
$ws.Range("B8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()

# Renumber the "م" column for the rows that followed (they kept their data,
# just shifted down by one row, so only the running index changes).
for ($i = 9; $i -le 18; $i++) {
    $ws.Cells.Item($i, 1).Value2 = $i - 3
}

# Fix the totals row: recompute the sum of column L across the data rows.
$total = 0
for ($i = 4; $i -le 18; $i++) {
    $v = $ws.Cells.Item($i, 12).Value2
    if ($v -ne $null -and $v -ne "") {
        $total = $total + $v
    }
}
$ws.Cells.Item(19, 11).Value2 = $total

# Excel re-wraps the footer text once the sheet gains a row, shrinking its
# row height slightly; match that final layout.
$ws.Rows.Item(20).RowHeight = 16.5
